# Generate Report for Handoff
# Update the localization-status workbook to reflect that the
# b0257c55-... file has now been handed off (zh-cn and de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-15 20:11:21"
# Columns E (zh-cn) and F (de-de) widen to fit the new "Ready for handoff" text
$wsOverview.Columns.Item(5).ColumnWidth = 16.37
$wsOverview.Columns.Item(6).ColumnWidth = 16.37

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-15 20:11:17"
# Column C (Status) widens to fit the new "Ready for handoff" text
$wsZhCn.Columns.Item(3).ColumnWidth = 16.37

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-15 20:11:21"
# Column C (Status) widens to fit the new "Ready for handoff" text
$wsDeDe.Columns.Item(3).ColumnWidth = 16.37
